$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 19330051920369
$ws.Range("B2").Value = "DE LOS SANTOS"
$ws.Range("C2").Value = "XOTLANIHUA"
$ws.Range("D2").Value = "JENNIFER"
$ws.Range("E2").Value = "ECOLOGÍA"
$ws.Range("F2").Value = "4ARHV"
$ws.Range("G2").Value = 1
